$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 42-45 describing get/set DSN & MODULE_SN scripts ---

# Row 42: get_DSN.sh
$ws.Range("B42").Value = "get_DSN.sh"
$ws.Range("C42").Value = "./get_DSN.sh"
$ws.Range("D42").Value = "adb shell /etc/factory-test/lvp15/get_DSN.sh"
$ws.Range("E42").Value = "get Device serial number"
$ws.Range("E42").Interior.Color = 65535

# Row 43: get_MODULE_SN.sh
$ws.Range("B43").Value = "get_MODULE_SN.sh"
$ws.Range("C43").Value = "./get_MODULE_SN.sh"
$ws.Range("D43").Value = "adb shell /etc/factory-test/lvp15/get_MODULE_SN.sh"
$ws.Range("E43").Value = "get Moduel serial number"
$ws.Range("E43").Interior.Color = 65535

# Row 44: set_DSN.sh
$ws.Range("B44").Value = "set_DSN.sh"
$ws.Range("C44").Value = "./set_DSN.sh"
$ws.Range("D44").Value = "adb shell /etc/factory-test/lvp15/set_DSN.sh 1234"
$ws.Range("E44").Value = "if success return OK"
$ws.Range("E44").Interior.Color = 65535

# Row 45: set_MODULE_SN.sh
$ws.Range("B45").Value = "set_MODULE_SN.sh"
$ws.Range("C45").Value = "./set_MODULE_SN.sh"
$ws.Range("D45").Value = "adb shell /etc/factory-test/lvp15/set_MODULE_SN.sh ZTQ03M130144"
$ws.Range("E45").Value = "if success return OK"
$ws.Range("E45").Interior.Color = 65535

# --- Mark obsolete / modified command rows with strikethrough formatting ---

# Rows 19-20 (no pre-existing fill) -> strikethrough only
$ws.Range("A19:D19").Font.Strikethrough = $true
$ws.Range("A20:D20").Font.Strikethrough = $true

# Rows 33-34: column A has no fill (strikethrough only),
# columns B:E already highlighted yellow (strikethrough + keep fill)
$ws.Range("A33").Font.Strikethrough = $true
$ws.Range("B33:E33").Font.Strikethrough = $true
$ws.Range("A34").Font.Strikethrough = $true
$ws.Range("B34:E34").Font.Strikethrough = $true

# Rows 40-41: columns A:D have no fill (strikethrough only),
# column E already highlighted yellow (strikethrough + keep fill)
$ws.Range("A40:D40").Font.Strikethrough = $true
$ws.Range("E40").Font.Strikethrough = $true
$ws.Range("A41:D41").Font.Strikethrough = $true
$ws.Range("E41").Font.Strikethrough = $true

# --- Update the sheet view so the newly added rows are visible ---
$null = $ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("E47").Select()
